# Update the two registrant e-mail addresses in the data list.
# (MongoDB-backed form data refresh — only the "email" column values change;
#  everything else on the sheet — names, phone numbers, hyperlink targets,
#  styles — stays exactly as it was.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (Abhay Kumar) first, then row 2 (Nitya Ranjan) — preserves the same
# shared-string insertion order produced by the original edit.
$ws.Range("C3").Value = "abhay43@gmail.com"
$ws.Range("C2").Value = "nityaranjan163@gmail.com"
